# Apply the edit described by the diff:
# - Add new column BB, mirroring column BA for rows 2-81 (same values)
# - Set a new header date in BB1 (copy BA1's format, then set new date)
# - Update BB82 and BB83 with new (non-mirrored) values (copy BA's format)
# - Add a brand-new row 84 with a date in A84 (copy A83's format) and a
#   value in BB84 (copy BA83's format)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header date for column BB (serial 45986 -> 2025-11-25).
# Copy BA1 first so BB1 picks up the same style (date format + border/bold).
$ws.Cells.Item(1, 53).Copy($ws.Cells.Item(1, 54))
$ws.Cells.Item(1, 54).Value = 45986

# Mirror column BA into column BB for rows 2 through 81 (inclusive).
# Copy() brings along BA's cell formatting (none, in this case) together
# with the value, so this reproduces both value and style faithfully.
for ($r = 2; $r -le 81; $r++) {
    $ws.Cells.Item($r, 53).Copy($ws.Cells.Item($r, 54))
}

# Row 82 and 83 get their own distinct BB values (not mirrored from BA).
# Still copy from BA to inherit formatting, then overwrite with the new value.
$ws.Cells.Item(82, 53).Copy($ws.Cells.Item(82, 54))
$ws.Cells.Item(82, 54).Value = 0.1765865160815849

$ws.Cells.Item(83, 53).Copy($ws.Cells.Item(83, 54))
$ws.Cells.Item(83, 54).Value = 0.2412052862208469

# New row 84: date in column A (style copied from A83), value in column BB
# only (style copied from BA83, which has no explicit style).
$ws.Cells.Item(83, 1).Copy($ws.Cells.Item(84, 1))
$ws.Cells.Item(84, 1).Value = 45884

$ws.Cells.Item(83, 53).Copy($ws.Cells.Item(84, 54))
$ws.Cells.Item(84, 54).Value = 0.768168485846715
